# "fichier excel non remplir" — clear the pre-filled solver results on the
# "Résultats" sheet: the optimal-benefit cell (B1) and the whole allocation
# matrix (B3:K12) should be empty again, leaving only the row/column labels
# (A3:A12, "Type 1".."Type 10") and the header row (row 2) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Résultats")

$ws.Range("B1").ClearContents()
$ws.Range("B3:K12").ClearContents()
